$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.262.40'
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").Value = '3.495.67'
$ws.Range("E3").Value = '  -0.61%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''590.63'
$ws.Range("E5").Value = '  +0.82%  '
$ws.Range("D6").Value = '''133.61'
$ws.Range("E6").Value = '  -0.74%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -0.47%  '
$ws.Range("D9").Value = '''7.34'
$ws.Range("E9").Value = '  +3.36%  '
$ws.Range("E10").Value = '  -0.11%  '
$ws.Range("E11").Value = '  +2.33%  '
$ws.Range("E13").Value = '  +1.15%  '
$ws.Range("E14").Value = '  +0.78%  '
$ws.Range("D15").Value = '3.496.24'
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("D16").Value = '64.343.70'
$ws.Range("D17").Value = '''25.64'
$ws.Range("E17").Value = '  -6.61%  '
$ws.Range("D18").Value = '''9.84'
$ws.Range("D19").Value = '''5.75'
$ws.Range("E19").Value = '  +2.44%  '
$ws.Range("D20").Value = '''13.48'
$ws.Range("E20").Value = '  -2.76%  '
$ws.Range("D21").Value = '''392.29'
$ws.Range("E21").Value = '  +2.43%  '
$ws.Range("D22").Value = '''0.571'
$ws.Range("E22").Value = '  +0.49%  '
$ws.Range("D23").Value = '3.634.79'
$ws.Range("E23").Value = '  -0.71%  '
$ws.Range("E24").Value = '  +0.61%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("E28").Value = '  -1.74%  '
$ws.Range("E29").Value = '  +1.52%  '
$ws.Range("E30").Value = '  -2.23%  '
$ws.Range("E31").Value = '  -7.30%  '
$ws.Range("D32").Value = '3.515.90'
$ws.Range("E32").Value = '  -0.50%  '
$ws.Range("E33").Value = '  +5.36%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").Value = '''23.38'
$ws.Range("E35").Value = '  -0.71%  '
$ws.Range("D36").Value = '''5.17'
$ws.Range("E36").Value = '  -4.29%  '
$ws.Range("D37").Value = '''6.86'
$ws.Range("E37").Value = '  -0.98%  '
$ws.Range("E38").Value = '  -0.93%  '
$ws.Range("D39").Value = '''166.66'
$ws.Range("E39").Value = '  +4.12%  '
$ws.Range("D40").Value = '''0.0780'
$ws.Range("E40").Value = '  -0.97%  '
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("D43").Value = '''25.11'
$ws.Range("E43").Value = '  -5.64%  '
$ws.Range("D44").Value = '''4.38'
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("E45").Value = '  +3.35%  '
$ws.Range("E46").Value = '  -3.49%  '
$ws.Range("E47").Value = '  -0.67%  '
$ws.Range("D48").Value = '2.368.87'
$ws.Range("E48").Value = '  -4.56%  '
$ws.Range("E49").Value = '  -2.95%  '
$ws.Range("E50").Value = '  -1.76%  '
$ws.Range("E51").Value = '  -1.37%  '
